# Natmi following Dr Hou advice
# Recompute the Gdnf-Gfra1 ligand-receptor table for the new cluster set
# (ECs / sCs / FAPs). Row 2 keeps its "Sending=ECs, Target=ECs" identity but
# gets fresh numbers, and rows 3-7 are new Sending x Target combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("ECs", "Gdnf", "Gfra1", "ECs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 2, 0.6666666666666666, 0.3392376666666667, 1.017713, 0.01394884456468018, 0.01394884456468017, 0.0397207729938889, 0.357486956945, 0.00130450201756229, 0.00130450201756229)
  3 = @("ECs", "Gdnf", "Gfra1", "FAPs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 3, 1, 20.328499, 60.985497, 0.8358714277529804, 0.8358714277529804, 2.380230067078333, 21.422070603705, 0.07817105989462547, 0.07817105989462546)
  4 = @("ECs", "Gdnf", "Gfra1", "sCs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 3, 1, 3.65239, 10.95717, 0.1501797276823394, 0.1501797276823394, 0.4276522577833333, 3.84887032005, 0.01404487352698943, 0.01404487352698943)
  5 = @("sCs", "Gdnf", "Gfra1", "ECs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 2, 0.6666666666666666, 0.3392376666666667, 1.017713, 0.01394884456468018, 0.01394884456468017, 0.3850074995741111, 3.465067496167, 0.01264434254711789, 0.01264434254711788)
  6 = @("sCs", "Gdnf", "Gfra1", "FAPs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 3, 1, 20.328499, 60.985497, 0.8358714277529804, 0.8358714277529804, 23.07121330891366, 207.640919780223, 0.757700367858355, 0.7577003678583549)
  7 = @("sCs", "Gdnf", "Gfra1", "sCs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 3, 1, 3.65239, 10.95717, 0.1501797276823394, 0.1501797276823394, 4.145169241336666, 37.30652317202999, 0.13613485415535, 0.1361348541553499)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item([int]$r, $i + 1).Value = $vals[$i]
  }
}
